$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Result")

$val_C2 = @'
10737510
'@
$val_C3 = @'
10737514
'@
$val_C4 = @'
10737523
'@
$val_E4 = @'
PASS
'@
$val_F4 = @'
Cannot invoke "org.openqa.selenium.WebElement.getText()" because the return value of "connect_OrderProcessNonSPL.P3P.isElementPresent(String)" is null
'@
$val_C5 = @'
10737551
'@
$val_F5 = @'
Cannot invoke "org.openqa.selenium.WebElement.getText()" because the return value of "connect_OCBaseMethods.TCAcknowledge.isElementPresent(String)" is null
'@
$val_C11 = @'
10737605
'@
$val_F11 = @'
no such element: Unable to locate element: {"method":"css selector","selector":"#lblServiceID"}
  (Session info: chrome=124.0.6367.158)
For documentation on this error, please visit: https://www.seleniumhq.org/exceptions/no_such_element.html
Build info: version: '3.141.59', revision: 'e82be7d358', time: '2018-11-14T08:17:03'
System info: host: 'SIPL92', ip: '10.212.140.10', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '20'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 124.0.6367.158, chrome: {chromedriverVersion: 124.0.6367.207 (a9001a6e39f..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:49723}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 7196eebf227b25c9aa98a113859bc596
*** Element info: {Using=id, value=lblServiceID}
'@
$val_C12 = @'
10737568
'@
$val_F12 = @'
no such element: Unable to locate element: {"method":"css selector","selector":"#lblServiceID"}
  (Session info: chrome=124.0.6367.158)
For documentation on this error, please visit: https://www.seleniumhq.org/exceptions/no_such_element.html
Build info: version: '3.141.59', revision: 'e82be7d358', time: '2018-11-14T08:17:03'
System info: host: 'SIPL92', ip: '10.212.140.10', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '20'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 124.0.6367.158, chrome: {chromedriverVersion: 124.0.6367.207 (a9001a6e39f..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:49723}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 7196eebf227b25c9aa98a113859bc596
*** Element info: {Using=id, value=lblServiceID}
'@
$val_C13 = @'
10737579
'@
$val_F13 = @'
no such element: Unable to locate element: {"method":"css selector","selector":"#lblServiceID"}
  (Session info: chrome=124.0.6367.158)
For documentation on this error, please visit: https://www.seleniumhq.org/exceptions/no_such_element.html
Build info: version: '3.141.59', revision: 'e82be7d358', time: '2018-11-14T08:17:03'
System info: host: 'SIPL92', ip: '10.212.140.10', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '20'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 124.0.6367.158, chrome: {chromedriverVersion: 124.0.6367.207 (a9001a6e39f..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:49723}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 7196eebf227b25c9aa98a113859bc596
*** Element info: {Using=id, value=lblServiceID}
'@
$val_C14 = @'
10737592
'@
$val_F14 = @'
no such element: Unable to locate element: {"method":"css selector","selector":"#lblServiceID"}
  (Session info: chrome=124.0.6367.158)
For documentation on this error, please visit: https://www.seleniumhq.org/exceptions/no_such_element.html
Build info: version: '3.141.59', revision: 'e82be7d358', time: '2018-11-14T08:17:03'
System info: host: 'SIPL92', ip: '10.212.140.10', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '20'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 124.0.6367.158, chrome: {chromedriverVersion: 124.0.6367.207 (a9001a6e39f..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:49723}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 7196eebf227b25c9aa98a113859bc596
*** Element info: {Using=id, value=lblServiceID}
'@
$val_E18 = @'
FAIL
'@
$val_C24 = @'
142049019
'@
$val_E24 = @'
FAIL
'@
$val_F24 = @'
Expected condition failed: waiting for visibility of all elements located by By.xpath: //*[@class="OpenCloseClass dropdown open"]//ul (tried for 60 second(s) with 500 milliseconds interval)
Build info: version: '3.141.59', revision: 'e82be7d358', time: '2018-11-14T08:17:03'
System info: host: 'SIPL92', ip: '10.212.140.10', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '20'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 124.0.6367.158, chrome: {chromedriverVersion: 124.0.6367.207 (a9001a6e39f..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:54212}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 9e124d070523c83c50bc6af0113b0bdc
'@
$val_F26 = @'
Expected condition failed: waiting for visibility of [[ChromeDriver: chrome on WINDOWS (27b24f416550d1afef5f0b758a15a28f)] -> id: GreyTick] (tried for 60 second(s) with 500 milliseconds interval)
Build info: version: '3.141.59', revision: 'e82be7d358', time: '2018-11-14T08:17:03'
System info: host: 'SIPL92', ip: '10.212.140.10', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '20'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 124.0.6367.158, chrome: {chromedriverVersion: 124.0.6367.207 (a9001a6e39f..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:50254}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 27b24f416550d1afef5f0b758a15a28f
'@

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = $val_C2
$ws.Range("C2").Style = "Normal"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = $val_C3
$ws.Range("C3").Style = "Normal"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = $val_C4
$ws.Range("C4").Style = "Normal"
$ws.Range("E4").Value = $val_E4
$ws.Range("F4").Value = $val_F4
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = $val_C5
$ws.Range("C5").Style = "Normal"
$ws.Range("F5").Value = $val_F5
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = $val_C11
$ws.Range("C11").Style = "Normal"
$ws.Range("F11").Value = $val_F11
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = $val_C12
$ws.Range("C12").Style = "Normal"
$ws.Range("F12").Value = $val_F12
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = $val_C13
$ws.Range("C13").Style = "Normal"
$ws.Range("F13").Value = $val_F13
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = $val_C14
$ws.Range("C14").Style = "Normal"
$ws.Range("F14").Value = $val_F14
$ws.Range("E18").Value = $val_E18
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = $val_C24
$ws.Range("C24").Style = "Normal"
$ws.Range("E24").Value = $val_E24
$ws.Range("F24").Value = $val_F24
$ws.Range("F26").Value = $val_F26
